$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the value for E2 (new shared string "This is a Git Demo")
$ws.Range("E2").Value = "This is a Git Demo"

# Move the selection to E3, as captured in the sheetView after the edit
$ws.Range("E3").Select()
